$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Global font rename: TimesNewToman -> Times New Roman (every run in body)
# ---------------------------------------------------------------------------
$fontFind = $d.Content.Find
$fontFind.ClearFormatting()
$fontFind.Font.Name = "TimesNewToman"
$fontFind.Replacement.ClearFormatting()
$fontFind.Replacement.Font.Name = "Times New Roman"
$fontFind.Execute("", $false, $false, $false, $false, $false, $true, 1, $true, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Title paragraph
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(1).Range.Text = 'A Journey Through Molecular Medicine: Connecting Health and Chemistry'

# ---------------------------------------------------------------------------
# 3) Author name paragraph
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(2).Range.Text = 'Dr. Eleanor Alvarez'

# ---------------------------------------------------------------------------
# 4) Author email paragraph
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(3).Range.Text = 'eleanor.alvarez@wright.edu'

# ---------------------------------------------------------------------------
# 5) Body paragraph (long, multi-sentence, with two double line-breaks)
# ---------------------------------------------------------------------------
$bodyText = 'In the tapestry of human existence, health occupies a central position, entwined with intricate biological processes. Understanding the molecular underpinnings of life offers a key to unravelling the mysteries of disease and paving the path towards novel treatments. Molecular medicine, an intersection of chemistry and biology, empowers us to decipher this symphony of life, illuminating the molecular mechanisms that govern health and illness. As we delve into the realm of molecular medicine, we embark on a journey filled with discoveries, challenges, and boundless opportunities.' + [char]11 + '' + [char]11 + 'At the heart of molecular medicine lies the chemistry of life. Molecules, the building blocks of life, engage in intricate interactions, dictating the symphony of biological processes. From DNA''s genetic code to the intricate machinery of proteins, chemistry offers a language to decode the language of life. Molecular medicine empowers us to manipulate these molecules, precisely targeting them to combat disease and promote health. By understanding the molecular mechanisms of disease, we can unravel the enigma of illness, paving the way for effective therapies and interventions.' + [char]11 + '' + [char]11 + 'Furthermore, molecular medicine offers a unique lens through which we can view and comprehend human health. It enables us to delve into the molecular interactions of the body, deciphering the intricate dance of cells, tissues, and organs. By examining the molecular basis of disease, we gain insights into the complex interplay of genetic, environmental, and lifestyle factors that influence our well-being. This understanding empowers us to intervene at the molecular level, preventing and treating diseases with greater precision and efficacy.'
$d.Paragraphs.Item(5).Range.Text = $bodyText

# ---------------------------------------------------------------------------
# 6) "Summary" heading paragraph - unchanged text, already correct
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 7) Summary body paragraph
# ---------------------------------------------------------------------------
$summaryText = 'Molecular medicine stands as a beacon of hope in the quest for understanding and treating diseases. By harnessing the power of chemistry and biology, molecular medicine offers a deeper understanding of the molecular basis of life and illness. It enables us to manipulate molecules, precisely targeting them to combat disease and promote health. Furthermore, molecular medicine provides a unique perspective on human health, allowing us to examine the intricate interactions of the body''s molecular machinery. With molecular medicine as our guide, we embark on a transformative journey towards better health and well-being for humanity.'
$d.Paragraphs.Item(7).Range.Text = $summaryText

# ---------------------------------------------------------------------------
# 8) Append a new, fully empty trailing paragraph after the summary text
#    (use ^p via Find/Replace so no run is left behind in the new paragraph)
# ---------------------------------------------------------------------------
$endFind = $d.Paragraphs.Item(7).Range.Find
$endFind.ClearFormatting()
$endFind.Replacement.ClearFormatting()
$endFind.Execute("well-being for humanity.", $false, $false, $false, $false, $false, $true, 1, $false, "well-being for humanity.^p", 2) | Out-Null
